$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Header rename: "FILTER/UNIT" -> "FILTERS/UNITS" (D1)
# ---------------------------------------------------------------------------
$ws.Range("D1").Value = "FILTERS/UNITS"

# Column D widened slightly to fit the longer header text.
$ws.Columns("D").ColumnWidth = 13.7

# ---------------------------------------------------------------------------
# 2. New model block "Model 1.12" (rows 95-100), same 6-row shape as the
#    preceding block (rows 88-93), so copy that block's formatting first.
# ---------------------------------------------------------------------------
$ws.Range("A88:M93").Copy()
$ws.Range("A95:M100").PasteSpecial(-4122)
# Re-apply the exact bottom-border formatting from row 93 onto row 100, since
# that row needs every column populated (it didn't get it from the 6-row
# block copy in a couple of spots it already matches, but doing it again is
# harmless and guarantees an exact match).
$ws.Range("A93:M93").Copy()
$ws.Range("A100:M100").PasteSpecial(-4122)

$ws.Range("A95").Value = "Model 1.12"
$ws.Range("B95").Value = "(64,64,1)"
$ws.Range("C95").Value = 16

$ws.Range("A96").Value = "Conv 1"
$ws.Range("D96").Value = 64
$ws.Range("E96").Value = "(3,3)"
$ws.Range("F96").Value = "(2,2)"
$ws.Range("G96").Value = "no"
$ws.Range("H96").Value = 0.2

$ws.Range("A97").Value = "Conv 2"
$ws.Range("D97").Value = 64
$ws.Range("E97").Value = "(3,3)"
$ws.Range("F97").Value = "(2,2)"
$ws.Range("G97").Value = "no"
$ws.Range("H97").Value = 0.3

$ws.Range("A98").Value = "Conv 3"
$ws.Range("D98").Value = 64
$ws.Range("E98").Value = "(5,5)"
$ws.Range("F98").Value = "(2,2)"
$ws.Range("G98").Value = "no"
$ws.Range("H98").Value = 0.5

$ws.Range("A99").Value = "Dense"
$ws.Range("D99").Value = 128
$ws.Range("H99").Value = 0.5

$ws.Range("A100").Value = "Output"
$ws.Range("I100").Value = 0.2225
$ws.Range("J100").Value = 0.9072
$ws.Range("K100").Value = 0.3057
$ws.Range("L100").Value = 0.8699
$ws.Range("M100").Value = 15

# ---------------------------------------------------------------------------
# 3. New model block "Model 1.13" (rows 102-107).
# ---------------------------------------------------------------------------
$ws.Range("A88:M93").Copy()
$ws.Range("A102:M107").PasteSpecial(-4122)
$ws.Range("A93:M93").Copy()
$ws.Range("A107:M107").PasteSpecial(-4122)

$ws.Range("A102").Value = "Model 1.13"
$ws.Range("B102").Value = "(128, 128, 1)"
$ws.Range("C102").Value = 16

$ws.Range("A103").Value = "Conv 1"
$ws.Range("D103").Value = 16
$ws.Range("E103").Value = "(3,3)"
$ws.Range("F103").Value = "(2,2)"
$ws.Range("G103").Value = "no"
$ws.Range("H103").Value = 0.2

$ws.Range("A104").Value = "Conv 2"
$ws.Range("D104").Value = 32
$ws.Range("E104").Value = "(3,3)"
$ws.Range("F104").Value = "(2,2)"
$ws.Range("G104").Value = "no"
$ws.Range("H104").Value = 0.2

$ws.Range("A105").Value = "Conv 3"
$ws.Range("D105").Value = 64
$ws.Range("E105").Value = "(5,5)"
$ws.Range("F105").Value = "(2,2)"
$ws.Range("G105").Value = "no"
$ws.Range("H105").Value = 0.3

$ws.Range("A106").Value = "Dense"
$ws.Range("D106").Value = 256
$ws.Range("H106").Value = 0.5

$ws.Range("A107").Value = "Output"
$ws.Range("I107").Value = 0.0902
$ws.Range("J107").Value = 0.9692
$ws.Range("K107").Value = 0.3538
$ws.Range("L107").Value = 0.8682
$ws.Range("M107").Value = 15

# ---------------------------------------------------------------------------
# 4. View state: keep the frozen header, move the active selection to the
#    new last cell, matching where the author ended up after the edit.
# ---------------------------------------------------------------------------
$ws.Range("M107").Select()
